# update document ID numbers
#
# 1) The "updated" date field (datetimeFigureOut) shown on every slide's
#    footer is sourced from the Slide Master and each of its Custom
#    Layouts. Bump it from 11/30/20 -> 1/12/21 everywhere it appears.
# 2) Several placeholder "xxxx"-style document IDs on the dashboard slide
#    get filled in with their real numbers.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "11/30/20") {
                $tr.Text = "1/12/21"
            }
        }
    }
}

# Slide Master footer date field.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Custom Layout has its own copy of the date field.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Dashboard slide: fill in the real document ID numbers.
$s = $p.Slides.Item(1)

# LOGOS Software Requirements Specification (SRS)
$sh = $s.Shapes.Item("Rectangle 23")
$found = $sh.TextFrame.TextRange.Find("SPC-xxxx")
$found.Text = "SPC-2979"

# LOGOS Software Design Description (SDD)
$sh = $s.Shapes.Item("Rectangle 27")
$found = $sh.TextFrame.TextRange.Find("SDD-xxx")
$found.Text = "SDD-559"

# LOGOS User Documentation
$sh = $s.Shapes.Item("Rectangle 44")
$found = $sh.TextFrame.TextRange.Find("INL/EXT-xx-xxxxx ")
$found.Text = "INL/EXT-20-61001"

# LOGOS Requirements Traceability Matrix (RTM)
$sh = $s.Shapes.Item("Rectangle 81")
$origHeight = $sh.Height
$found = $sh.TextFrame.TextRange.Find("SPC-xxxx")
$found.Text = "SPC-2979"
# This shape's text box has auto-fit enabled, and merging the two runs
# that made up "SPC-xxxx" into the single run "SPC-2979" makes the
# runtime recompute the auto-fit height. Put the box back to its
# original size since the real edit never touched the shape geometry.
$sh.Height = $origHeight

# LOGOS Configuration Items List
$sh = $s.Shapes.Item("Rectangle 66")
$found = $sh.TextFrame.TextRange.Find("LST-xxxx")
$found.Text = "LST-1291"
